$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Status")

# --- Insert the first new point: "ADDED POINT FOR 115KV LINE POS." ---
# It goes in right before the existing "IINYO 115KV CB" row (currently row 4).
$target1 = $ws.Range("B1:B200").Find("IINYO 115KV CB")
$row1 = $target1.Row
$ws.Rows.Item($row1).Insert()
$ws.Range("B" + $row1).Value = "ADDED POINT FOR 115KV LINE POS."

# --- Insert the second new point: "WEST 115KV BUS DIFF SYS B RLY FAIL" ---
# It goes in right before the existing "WEST 115KV BUS DIFF SYS A RLY TRIP" row.
$target2 = $ws.Range("B1:B200").Find("WEST 115KV BUS DIFF SYS A RLY TRIP")
$row2 = $target2.Row
$ws.Rows.Item($row2).Insert()
$ws.Range("B" + $row2).Value = "WEST 115KV BUS DIFF SYS B RLY FAIL"

# --- Renumber the "Point Number" column (A) sequentially for all data rows ---
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row   # xlUp = -4162
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}
